$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.503.34"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.871.05"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -1.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5068"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3913"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08338"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.37"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.104"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.187"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "1.870.30"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.30"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.244"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.50"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001098"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06717"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.62"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.913"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("D23").Value = "28.529.15"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.192"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.04%  "
$ws.Range("D26").Value = "2.076.06"
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.63"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.418"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.45"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.038"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.772"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.625"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02439"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06547"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.959"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2158"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.025"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.183"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.240"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6357"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.07"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5982"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.96"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.676"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.999"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.209"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.03"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.154"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.25%  "
